# Auto-generated edit script: updates cached market-data values
# in the per-profession Leve profit tables to match the scheduled
# data-refresh commit.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 1445.7667
$ws.Range("I137").Value = 882.75
$ws.Range("J137").Value = 1821.1111
$ws.Range("K137").Value = 2648.25
$ws.Range("L137").Value = 5463.3333
$ws.Range("M137").Value = -98.25
$ws.Range("N137").Value = -10563.3333

# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 1820.1
$ws.Range("I138").Value = 644.1818
$ws.Range("J138").Value = 2399.2837
$ws.Range("K138").Value = 1932.5454
$ws.Range("L138").Value = 7197.8511
$ws.Range("M138").Value = 3207.4546
$ws.Range("N138").Value = -17477.8511

$ws = $wb.Worksheets.Item("ARM")
# Row 97 (Leve Item ID 19941)
$ws.Range("H97").Value = 14493408
$ws.Range("I97").Value = 22222884
$ws.Range("J97").Value = 640.25
$ws.Range("K97").Value = 22222884
$ws.Range("L97").Value = 640.25
$ws.Range("M97").Value = -22222388
$ws.Range("N97").Value = -1632.25

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 1442.8077
$ws.Range("I31").Value = 1068.6
$ws.Range("J31").Value = 1747.3954
$ws.Range("K31").Value = 1068.6
$ws.Range("L31").Value = 1747.3954
$ws.Range("M31").Value = -773.5999999999999
$ws.Range("N31").Value = -2337.3954

# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 1442.8077
$ws.Range("I34").Value = 1068.6
$ws.Range("J34").Value = 1747.3954
$ws.Range("K34").Value = 1068.6
$ws.Range("L34").Value = 1747.3954
$ws.Range("M34").Value = -866.5999999999999
$ws.Range("N34").Value = -2151.3954

# Row 107 (Leve Item ID 27689)
$ws.Range("H107").Value = 594.7273
$ws.Range("I107").Value = 501.46155
$ws.Range("K107").Value = 501.46155
$ws.Range("M107").Value = 1418.53845

$ws = $wb.Worksheets.Item("CUL")
# Row 22 (Leve Item ID 4697)
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()

# Row 23 (Leve Item ID 4858)
$ws.Range("H23").Value = 1160.9166
$ws.Range("I23").Value = 1938.6666
$ws.Range("J23").Value = 383.16666
$ws.Range("K23").Value = 5815.9998
$ws.Range("L23").Value = 1149.49998
$ws.Range("M23").Value = -5580.9998
$ws.Range("N23").Value = -1619.49998

# Row 26 (Leve Item ID 4746)
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()

# Row 27 (Leve Item ID 4697)
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("N27").ClearContents()

# Row 32 (Leve Item ID 4731)
$ws.Range("H32").Value = 534.3333
$ws.Range("J32").Value = 1003
$ws.Range("L32").Value = 3009
$ws.Range("N32").Value = -3575

# Row 33 (Leve Item ID 4867)
$ws.Range("H33").Value = 151
$ws.Range("I33").Value = 179.57143
$ws.Range("K33").Value = 1077.42858
$ws.Range("M33").Value = -794.42858

# Row 38 (Leve Item ID 4860)
$ws.Range("H38").Value = 135.12
$ws.Range("J38").Value = 100
$ws.Range("L38").Value = 300
$ws.Range("N38").Value = -994

# Row 39 (Leve Item ID 4712)
$ws.Range("H39").Value = 8895.097
$ws.Range("J39").Value = 8895.097
$ws.Range("L39").Value = 26685.291
$ws.Range("N39").Value = -27273.291

# Row 41 (Leve Item ID 4700)
$ws.Range("H41").Value = 642.2222
$ws.Range("I41").Value = 193.33333
$ws.Range("J41").Value = 866.6667
$ws.Range("K41").Value = 579.99999
$ws.Range("L41").Value = 2600.0001
$ws.Range("M41").Value = -241.99999
$ws.Range("N41").Value = -3276.0001

# Row 50 (Leve Item ID 4725)
$ws.Range("H50").Value = 454.44446
$ws.Range("I50").Value = 446.66666
$ws.Range("J50").Value = 458.33334
$ws.Range("K50").Value = 1339.99998
$ws.Range("L50").Value = 1375.00002
$ws.Range("M50").Value = -858.9999800000001
$ws.Range("N50").Value = -2337.00002

# Row 53 (Leve Item ID 4725)
$ws.Range("H53").Value = 454.44446
$ws.Range("I53").Value = 446.66666
$ws.Range("J53").Value = 458.33334
$ws.Range("K53").Value = 1339.99998
$ws.Range("L53").Value = 1375.00002
$ws.Range("M53").Value = -858.9999800000001
$ws.Range("N53").Value = -2337.00002

# Row 57 (Leve Item ID 4655)
$ws.Range("H57").Value = 4833.3335
$ws.Range("J57").Value = 4833.3335
$ws.Range("L57").Value = 14500.0005
$ws.Range("N57").Value = -15618.0005

# Row 96 (Leve Item ID 19816)
$ws.Range("H96").Value = 4200
$ws.Range("J96").Value = 4200
$ws.Range("L96").Value = 12600
$ws.Range("N96").Value = -16718

$ws = $wb.Worksheets.Item("GSM")
# Row 122 (Leve Item ID 36182)
$ws.Range("H122").Value = 1011127.75
$ws.Range("I122").Value = 1587989.2
$ws.Range("J122").Value = 1620
$ws.Range("K122").Value = 4763967.6
$ws.Range("L122").Value = 4860
$ws.Range("M122").Value = -4761517.6
$ws.Range("N122").Value = -9760

# Row 123 (Leve Item ID 34150)
$ws.Range("H123").Value = 10323
$ws.Range("J123").Value = 10323
$ws.Range("L123").Value = 10323
$ws.Range("N123").Value = -15223

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (Leve Item ID 36249)
$ws.Range("H7").Value = 3367.3684
$ws.Range("I7").Value = 2996.6667
$ws.Range("J7").Value = 3538.4614
$ws.Range("K7").Value = 2996.6667
$ws.Range("L7").Value = 3538.4614
$ws.Range("M7").Value = -2884.6667
$ws.Range("N7").Value = -3762.4614

# Row 70 (Leve Item ID 10811)
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

# Row 73 (Leve Item ID 10811)
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

# Row 126 (Leve Item ID 36249)
$ws.Range("H126").Value = 3367.3684
$ws.Range("I126").Value = 2996.6667
$ws.Range("J126").Value = 3538.4614
$ws.Range("K126").Value = 8990.000100000001
$ws.Range("L126").Value = 10615.3842
$ws.Range("M126").Value = -6520.000100000001
$ws.Range("N126").Value = -15555.3842

# Row 136 (Leve Item ID 44060)
$ws.Range("H136").Value = 3364.5945
$ws.Range("I136").Value = 1455.92
$ws.Range("J136").Value = 7341
$ws.Range("K136").Value = 4367.76
$ws.Range("L136").Value = 22023
$ws.Range("M136").Value = -1817.76
$ws.Range("N136").Value = -27123

$ws = $wb.Worksheets.Item("WVR")
# Row 126 (Leve Item ID 36210)
$ws.Range("H126").Value = 167883.5
$ws.Range("I126").Value = 201060.2
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 603180.6000000001
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -600710.6000000001
$ws.Range("N126").Value = -10940

# Row 128 (Leve Item ID 34563)
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 27780880
$ws.Range("I132").Value = 38464150
$ws.Range("K132").Value = 115392450
$ws.Range("M132").Value = -115389920
